$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 17 (shifts rows 17-28 down to 18-29),
# copying formatting from the row above (row 16, "IE")
$ws.Range("A17:C17").Insert(-4121, 0)  # xlShiftDown, xlFormatFromLeftOrAbove

# Set the new values for row 17: Region = "IS", B = 30.572, C = 32.65
# Leading apostrophe preserves the quotePrefix text style used by the other
# region-code cells in column A.
$ws.Cells.Item(17, 1).Value = "'IS"
$ws.Cells.Item(17, 2).Value = 30.571999999999999
$ws.Cells.Item(17, 3).Value = 32.65

$wb.Save()
